# Refresh the scraped "想去人数" (want-to-go count) figures in column F
# for the generated data, across the "展览" and "全部类型" sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibition) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 9487
$ws1.Range("F3").Value = 210
$ws1.Range("F4").Value = 25
$ws1.Range("F5").Value = 526
$ws1.Range("F6").Value = 464

# --- Sheet "全部类型" (all types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 9487
$ws4.Range("F3").Value = 210
$ws4.Range("F4").Value = 25
$ws4.Range("F5").Value = 526
$ws4.Range("F7").Value = 464
